$d = $word.ActiveDocument

# 1) Drop the stray empty run that trails "Template Transforms:" by forcing
#    a rewrite of that run (replacing the text with itself collapses the
#    paragraph down to a single populated run).
$d.Content.Find.Execute(
    "Template Transforms:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Template Transforms:",
    2) | Out-Null

# 2) Occurrence : Statement -> Occurrence : Transform (Mapping bullet).
$d.Content.Find.Execute(
    "(Context : Mapping, Occurrence : Statement, Attribute : Resource T, Value : Resource U);",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(Context : Mapping, Occurrence : Transform, Attribute : Resource T, Value : Resource U);",
    2) | Out-Null

# 3) Occurrence : Kind -> Occurrence : Mapping (Statement bullet).
$d.Content.Find.Execute(
    "(Context : Statement, Occurrence : Kind, Attribute : Resource, Value : Resource);",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(Context : Statement, Occurrence : Mapping, Attribute : Resource, Value : Resource);",
    2) | Out-Null
